$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from H1 into the new
# header cells I1 and J1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I ("I0") and J ("IF"), rows 2-24.
$data = @{
    2  = @(8, 8)
    3  = @(5, 7)
    4  = @(4, 5)
    5  = @(7, 8)
    6  = @(6, 6)
    7  = @(8, 8)
    8  = @(6, 7)
    9  = @(8, 8)
    10 = @(7, 7)
    11 = @(7, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(3, 4)
    17 = @(6, 7)
    18 = @(8, 10)
    19 = @(7, 7)
    20 = @(4, 5)
    21 = @(6, 8)
    22 = @(1, 5)
    23 = @(6, 9)
    24 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
